$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 554, shifting existing rows 554:613 down to 555:614.
$ws.Rows.Item(554).EntireRow.Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Cells.Item(554, 1).Value = 3
$ws.Cells.Item(554, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(554, 3).Value = "Coquimbo"
$ws.Cells.Item(554, 4).Value = 45194
$ws.Cells.Item(554, 5).Value = 5
$ws.Cells.Item(554, 6).Value = 100112012
$ws.Cells.Item(554, 7).Value = "Espinaca"
$ws.Cells.Item(554, 8).Value = "Sin especificar"
$ws.Cells.Item(554, 9).Value = "Primera"
$ws.Cells.Item(554, 10).Value = 90
$ws.Cells.Item(554, 11).Value = 4000
$ws.Cells.Item(554, 12).Value = 4500
$ws.Cells.Item(554, 13).Value = 4222
$ws.Cells.Item(554, 14).Value = "$/docena de atados (3 kilos)"
$ws.Cells.Item(554, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(554, 16).Value = 1407
$ws.Cells.Item(554, 17).Value = 3
$ws.Cells.Item(554, 18).Value = "Hortaliza"
